# Apply draft-mapping changes to StructureDefinition-ror-organization-closing-type.xlsx
#
# 1. Bump the "Date" metadata value on the Metadata sheet.
# 2. Add a new "Mapping: Spécification métier vers l'extension ROR ClosingType"
#    column (AL) to the Elements sheet, with a value only on the
#    Extension.value[x] row ("typeFermeture").

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the Date property -----------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-12T09:15:29+00:00"

# --- Elements sheet: add the new mapping column ----------------------------
$elements = $wb.Worksheets.Item("Elements")

# New column header in AL1 (column 38), and the one mapped value, on the
# Extension.value[x] row (row 6). AL2:AL5 are left blank, matching the
# other "Mapping" columns.
$elements.Range("AL1").Value = "Mapping: Spécification métier vers l'extension ROR ClosingType"
$elements.Range("AL6").Value = "typeFermeture"

# Match the formatting used by the rest of the table: header style for row
# 1, body style for the data rows.
$elements.Range("AK1").Copy()
$elements.Range("AL1").PasteSpecial(-4122)
$elements.Range("AK2:AK6").Copy()
$elements.Range("AL2:AL6").PasteSpecial(-4122)

# Approximate the "best fit" column width Excel would compute for the new
# column's content.
$elements.Columns.Item(38).AutoFit()
